# Auto-generated Word COM-interop script
# Updates the title date and all 100 arithmetic answers in the table

$d = $word.ActiveDocument

# --- Update the title paragraph (date line) ---
$titleRange = $d.Paragraphs.Item(1).Range
[void]$titleRange.MoveEnd(1, -1)
$titleRange.Text = "2025-11-02 Sunday"

# --- Update each answer cell in the single table (row-major order) ---
$newValues = @(
    "38-17=21",
    "59-50=9",
    "38-33=5",
    "36-27=9",
    "57-20=37",
    "87-78=9",
    "67+28=95",
    "63-22=41",
    "97-73=24",
    "49-41=8",
    "42+15=57",
    "75-34=41",
    "26+51=77",
    "81-2=79",
    "28+5=33",
    "0+5=5",
    "44-31=13",
    "14+67=81",
    "52+37=89",
    "42+23=65",
    "50-45=5",
    "10-1=9",
    "61-56=5",
    "27+41=68",
    "85-18=67",
    "80+1=81",
    "29+34=63",
    "36+49=85",
    "85-73=12",
    "5+21=26",
    "97-83=14",
    "71-11=60",
    "5+29=34",
    "24-10=14",
    "74+16=90",
    "1+3=4",
    "64-29=35",
    "44-7=37",
    "82-8=74",
    "28+8=36",
    "28+59=87",
    "44+0=44",
    "73+16=89",
    "39+41=80",
    "91-26=65",
    "4+57=61",
    "0+31=31",
    "63-10=53",
    "98-85=13",
    "20+66=86",
    "21-9=12",
    "11+8=19",
    "55+3=58",
    "57+42=99",
    "65+7=72",
    "7+73=80",
    "34-24=10",
    "27-16=11",
    "40-7=33",
    "56+37=93",
    "88+10=98",
    "61-34=27",
    "29+68=97",
    "15+57=72",
    "61+33=94",
    "84+3=87",
    "93-81=12",
    "4+81=85",
    "79-69=10",
    "27-24=3",
    "56-11=45",
    "10+12=22",
    "11+54=65",
    "70-47=23",
    "77-53=24",
    "87-45=42",
    "3+81=84",
    "53+12=65",
    "34+52=86",
    "35+13=48",
    "84-79=5",
    "87-86=1",
    "63-62=1",
    "1+87=88",
    "93-60=33",
    "34-17=17",
    "22+48=70",
    "76-41=35",
    "96-41=55",
    "59-28=31",
    "74+1=75",
    "72-59=13",
    "29+59=88",
    "67-43=24",
    "78-68=10",
    "7+37=44",
    "51-50=1",
    "24+14=38",
    "87-70=17",
    "86-57=29"
)

$t = $d.Tables.Item(1)
$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        [void]$rng.MoveEnd(1, -1)
        $rng.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}

Write-Output ("Updated " + $idx + " cells")
